$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet, matching the
# workbook.xml diff which appends "joint-scat" as sheetId 19 / rId19.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "joint-scat"

# Match the <outlinePr summaryBelow="1" summaryRight="1"/> used on every
# other sheet in the workbook.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# --- Header row + data rows (values) -----------------------------------
$ws.Cells.Item(1, 1).Value = "alg_code"
$ws.Cells.Item(1, 2).Value = "dog_bark"
$ws.Cells.Item(1, 3).Value = "footstep"
$ws.Cells.Item(1, 4).Value = "gunshot"
$ws.Cells.Item(1, 5).Value = "keyboard"
$ws.Cells.Item(1, 6).Value = "moving_motor_vehicle"
$ws.Cells.Item(1, 7).Value = "rain"
$ws.Cells.Item(1, 8).Value = "sneeze_cough"
$ws.Cells.Item(1, 9).Value = "avg_category_FAD"
$ws.Cells.Item(2, 1).Value = "Baseline"
$ws.Cells.Item(2, 2).Value = 6412.108822858143
$ws.Cells.Item(2, 3).Value = 16177.61467798351
$ws.Cells.Item(2, 4).Value = 226877.2591667986
$ws.Cells.Item(2, 5).Value = 39450.46834742428
$ws.Cells.Item(2, 6).Value = 417760.8205551166
$ws.Cells.Item(2, 7).Value = 1246.460676927177
$ws.Cells.Item(2, 8).Value = 58884.08675310281
$ws.Cells.Item(2, 9).Value = 109544.1170000301
$ws.Cells.Item(3, 1).Value = "TASys02"
$ws.Cells.Item(3, 2).Value = 2164.677705544587
$ws.Cells.Item(3, 3).Value = 11466.7227572848
$ws.Cells.Item(3, 4).Value = 12766.89027288463
$ws.Cells.Item(3, 5).Value = 120962.860210016
$ws.Cells.Item(3, 6).Value = 345588.3932784392
$ws.Cells.Item(3, 7).Value = 9572.154394896079
$ws.Cells.Item(3, 8).Value = 34048.72736894032
$ws.Cells.Item(3, 9).Value = 76652.91799828652
$ws.Cells.Item(4, 1).Value = "TASys03"
$ws.Cells.Item(4, 2).Value = 8326.955484781982
$ws.Cells.Item(4, 3).Value = 7925.430432487436
$ws.Cells.Item(4, 4).Value = 228332.4501374071
$ws.Cells.Item(4, 5).Value = 127955.5979676959
$ws.Cells.Item(4, 6).Value = 795967.2159420585
$ws.Cells.Item(4, 7).Value = 9027.467967408056
$ws.Cells.Item(4, 8).Value = 101063.4592286405
$ws.Cells.Item(4, 9).Value = 182656.9395943542
$ws.Cells.Item(5, 1).Value = "TASys08"
$ws.Cells.Item(5, 2).Value = 5475.915396630806
$ws.Cells.Item(5, 3).Value = 3494.91337050969
$ws.Cells.Item(5, 4).Value = 165748.3614948218
$ws.Cells.Item(5, 5).Value = 111226.6460575687
$ws.Cells.Item(5, 6).Value = 831672.3790306748
$ws.Cells.Item(5, 7).Value = 7507.48932138861
$ws.Cells.Item(5, 8).Value = 101363.203904979
$ws.Cells.Item(5, 9).Value = 175212.7012252247
$ws.Cells.Item(6, 1).Value = "TASys11"
$ws.Cells.Item(6, 2).Value = 7599.925545613903
$ws.Cells.Item(6, 3).Value = 150297.0603118276
$ws.Cells.Item(6, 4).Value = 297024.9416328784
$ws.Cells.Item(6, 5).Value = 24044.90400988286
$ws.Cells.Item(6, 6).Value = 837198.6016055684
$ws.Cells.Item(6, 7).Value = 8100.835438509881
$ws.Cells.Item(6, 8).Value = 4753.224838072987
$ws.Cells.Item(6, 9).Value = 189859.9276260506
$ws.Cells.Item(7, 1).Value = "TBSys09"
$ws.Cells.Item(7, 2).Value = 4490.060552633944
$ws.Cells.Item(7, 3).Value = 7256.478243895912
$ws.Cells.Item(7, 4).Value = 154775210.6032548
$ws.Cells.Item(7, 5).Value = 94578.96995438276
$ws.Cells.Item(7, 6).Value = 572297.5011757913
$ws.Cells.Item(7, 7).Value = 197.2951573284226
$ws.Cells.Item(7, 8).Value = 95317.25800536814
$ws.Cells.Item(7, 9).Value = 22221335.45233488
$ws.Cells.Item(8, 1).Value = "TBSys14"
$ws.Cells.Item(8, 2).Value = 5868.400010851554
$ws.Cells.Item(8, 3).Value = 255.1290328809373
$ws.Cells.Item(8, 4).Value = 3120145.402389309
$ws.Cells.Item(8, 5).Value = 105883.8738016963
$ws.Cells.Item(8, 6).Value = 761575.5037789284
$ws.Cells.Item(8, 7).Value = 2405.658471966781
$ws.Cells.Item(8, 8).Value = 66856.95550418066
$ws.Cells.Item(8, 9).Value = 580427.2747128305
$ws.Cells.Item(9, 1).Value = "TBSys18"
$ws.Cells.Item(9, 2).Value = 956.6505526141482
$ws.Cells.Item(9, 3).Value = 1110.401407768937
$ws.Cells.Item(9, 4).Value = 270358.8138694638
$ws.Cells.Item(9, 5).Value = 108498.1348472907
$ws.Cells.Item(9, 6).Value = 812192.4827031174
$ws.Cells.Item(9, 7).Value = 928.9746406285522
$ws.Cells.Item(9, 8).Value = 88940.92322370921
$ws.Cells.Item(9, 9).Value = 183283.7687492275
$ws.Cells.Item(10, 1).Value = "TBSys24"
$ws.Cells.Item(10, 2).Value = 3290.537615515859
$ws.Cells.Item(10, 3).Value = 2559.433608909461
$ws.Cells.Item(10, 4).Value = 129508.1255027545
$ws.Cells.Item(10, 5).Value = 100784.188118854
$ws.Cells.Item(10, 6).Value = 701252.0431087858
$ws.Cells.Item(10, 7).Value = 2720.61437615156
$ws.Cells.Item(10, 8).Value = 59086.40307576059
$ws.Cells.Item(10, 9).Value = 142743.0493438188

# --- Formatting: bold / bordered / centered header row & label column --
# (mirrors style index "1" used for row 1 and column A on every other sheet)
$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$labelRange = $ws.Range("A1:A10")
$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1

# --- Page margins (0.75/0.75/1/1 in, 0.5/0.5 in header/footer) ---------
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

$ws.Range("A1").Select() | Out-Null
